# Auto-generated edit script: apply scraped Leve profit-sheet updates across all job sheets
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 575.9286
$ws.Range("I38").Value = 274.15384
$ws.Range("K38").Value = 822.4615200000001
$ws.Range("M38").Value = -450.4615200000001
$ws.Range("H62").Value = 64291.438
$ws.Range("I62").Value = 84417.914
$ws.Range("J62").Value = 3912
$ws.Range("K62").Value = 84417.914
$ws.Range("L62").Value = 3912
$ws.Range("M62").Value = -83793.914
$ws.Range("N62").Value = -5160
$ws.Range("H65").Value = 64291.438
$ws.Range("I65").Value = 84417.914
$ws.Range("J65").Value = 3912
$ws.Range("K65").Value = 422089.57
$ws.Range("L65").Value = 19560
$ws.Range("M65").Value = -418969.57
$ws.Range("N65").Value = -25800
$ws.Range("H112").Value = 772113.75
$ws.Range("J112").Value = 836290.75
$ws.Range("L112").Value = 2508872.25
$ws.Range("N112").Value = -2511088.25
$ws.Range("H116").Value = 3243.6667
$ws.Range("I116").Value = 3481
$ws.Range("J116").Value = 3125
$ws.Range("K116").Value = 3481
$ws.Range("L116").Value = 3125
$ws.Range("M116").Value = -39
$ws.Range("N116").Value = -10009
$ws.Range("H124").Value = 104666.336
$ws.Range("J124").Value = 104666.336
$ws.Range("L124").Value = 104666.336
$ws.Range("N124").Value = -114486.336
$ws.Range("H132").Value = 3635.3225
$ws.Range("I132").Value = 3640.451
$ws.Range("K132").Value = 10921.353
$ws.Range("M132").Value = -8391.352999999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2841.152
$ws.Range("I32").Value = 3011.6511
$ws.Range("K32").Value = 3011.6511
$ws.Range("M32").Value = -2724.6511
$ws.Range("H42").Value = 40000
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H45").Value = 4780.8
$ws.Range("I45").Value = 2231.2856
$ws.Range("K45").Value = 2231.2856
$ws.Range("M45").Value = -1854.2856
$ws.Range("H97").Value = 1582.4706
$ws.Range("I97").Value = 1057.08
$ws.Range("J97").Value = 3041.889
$ws.Range("K97").Value = 1057.08
$ws.Range("L97").Value = 3041.889
$ws.Range("M97").Value = -561.0799999999999
$ws.Range("N97").Value = -4033.889
$ws.Range("H110").Value = 5709.846
$ws.Range("I110").Value = 2097.2856
$ws.Range("K110").Value = 2097.2856
$ws.Range("M110").Value = -52.28560000000016
$ws.Range("H122").Value = 3902.4285
$ws.Range("I122").Value = 4829.5
$ws.Range("J122").Value = 2666.3333
$ws.Range("K122").Value = 14488.5
$ws.Range("L122").Value = 7998.999899999999
$ws.Range("M122").Value = -12038.5
$ws.Range("N122").Value = -12898.9999
$ws.Range("H132").Value = 9155.0625
$ws.Range("I132").Value = 3771.3635
$ws.Range("K132").Value = 11314.0905
$ws.Range("M132").Value = -8784.0905

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2071.4546
$ws.Range("I20").Value = 1954.8823
$ws.Range("K20").Value = 1954.8823
$ws.Range("M20").Value = -1707.8823
$ws.Range("H86").Value = 2559.5
$ws.Range("I86").Value = 2496.375
$ws.Range("J86").Value = 2685.75
$ws.Range("K86").Value = 2496.375
$ws.Range("L86").Value = 2685.75
$ws.Range("M86").Value = -1373.375
$ws.Range("N86").Value = -4931.75
$ws.Range("H89").Value = 2559.5
$ws.Range("I89").Value = 2496.375
$ws.Range("J89").Value = 2685.75
$ws.Range("K89").Value = 12481.875
$ws.Range("L89").Value = 13428.75
$ws.Range("M89").Value = -6865.875
$ws.Range("N89").Value = -24660.75
$ws.Range("H92").Value = 74293.375
$ws.Range("J92").Value = 74293.375
$ws.Range("L92").Value = 74293.375
$ws.Range("N92").Value = -79285.375
$ws.Range("H94").Value = 1135.9546
$ws.Range("I94").Value = 1099.6316
$ws.Range("K94").Value = 1099.6316
$ws.Range("M94").Value = -648.6315999999999
$ws.Range("H122").Value = 73500
$ws.Range("J122").Value = 73500
$ws.Range("L122").Value = 73500
$ws.Range("N122").Value = -83300

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2799.7073
$ws.Range("I31").Value = 1223.0526
$ws.Range("J31").Value = 4161.364
$ws.Range("K31").Value = 1223.0526
$ws.Range("L31").Value = 4161.364
$ws.Range("M31").Value = -928.0526
$ws.Range("N31").Value = -4751.364
$ws.Range("H34").Value = 2799.7073
$ws.Range("I34").Value = 1223.0526
$ws.Range("J34").Value = 4161.364
$ws.Range("K34").Value = 1223.0526
$ws.Range("L34").Value = 4161.364
$ws.Range("M34").Value = -1021.0526
$ws.Range("N34").Value = -4565.364
$ws.Range("H102").Value = 96417.5
$ws.Range("J102").Value = 96417.5
$ws.Range("L102").Value = 96417.5
$ws.Range("N102").Value = -101285.5
$ws.Range("H131").Value = 91163
$ws.Range("J131").Value = 89326
$ws.Range("L131").Value = 89326
$ws.Range("N131").Value = -99406
$ws.Range("H132").Value = 1379.8
$ws.Range("I132").Value = 1535
$ws.Range("J132").Value = 1147
$ws.Range("K132").Value = 4605
$ws.Range("L132").Value = 3441
$ws.Range("M132").Value = -2075
$ws.Range("N132").Value = -8501

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 3678.1875
$ws.Range("J107").Value = 4477.923
$ws.Range("L107").Value = 13433.769
$ws.Range("N107").Value = -17273.769
$ws.Range("H132").Value = 2543.65
$ws.Range("I132").Value = 2003.5
$ws.Range("J132").Value = 2678.6875
$ws.Range("K132").Value = 18031.5
$ws.Range("L132").Value = 24108.1875
$ws.Range("M132").Value = -15501.5
$ws.Range("N132").Value = -29168.1875
$ws.Range("H134").Value = 2827.3572
$ws.Range("I134").Value = 2121.7693
$ws.Range("J134").Value = 12000
$ws.Range("K134").Value = 6365.3079
$ws.Range("L134").Value = 36000
$ws.Range("M134").Value = -1295.3079
$ws.Range("N134").Value = -46140

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 74166.164
$ws.Range("J136").Value = 74166.164
$ws.Range("L136").Value = 222498.492
$ws.Range("N136").Value = -227598.492

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 2399.5
$ws.Range("I13").Value = 2299
$ws.Range("K13").Value = 2299
$ws.Range("M13").Value = -2159
$ws.Range("H17").Value = 2169.3333
$ws.Range("I17").Value = 2254
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 2254
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = -2084
$ws.Range("N17").Value = -2340
$ws.Range("H43").Value = 28888
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 28888
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 28888
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -29274
$ws.Range("H74").Value = 25625
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 25909.092
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 25909.092
$ws.Range("M74").Value = -9002
$ws.Range("N74").Value = -27905.092
$ws.Range("H77").Value = 25625
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 25909.092
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 77727.276
$ws.Range("M77").Value = -25008
$ws.Range("N77").Value = -87711.276
$ws.Range("H136").Value = 5114.5557
$ws.Range("J136").Value = 16316.667
$ws.Range("L136").Value = 48950.001
$ws.Range("N136").Value = -54050.001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 83342570
$ws.Range("I2").Value = 10165.875
$ws.Range("J2").Value = 250007380
$ws.Range("K2").Value = 10165.875
$ws.Range("L2").Value = 250007380
$ws.Range("M2").Value = -10053.875
$ws.Range("N2").Value = -250007604
$ws.Range("H4").Value = 62503260
$ws.Range("I4").Value = 2036.3334
$ws.Range("J4").Value = 142861970
$ws.Range("K4").Value = 2036.3334
$ws.Range("L4").Value = 142861970
$ws.Range("M4").Value = -1923.3334
$ws.Range("N4").Value = -142862196
$ws.Range("H49").Value = 34444
$ws.Range("J49").Value = 34444
$ws.Range("L49").Value = 34444
$ws.Range("N49").Value = -34904
$ws.Range("H96").Value = 2844.16
$ws.Range("J96").Value = 3164.2856
$ws.Range("L96").Value = 3164.2856
$ws.Range("N96").Value = -5910.2856

